$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '35.079.91'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.73%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.857.51'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.17%  '

$ws.Range("E4").Value = '  -0.21%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.87'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.77%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.623'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.23%  '

$ws.Range("E7").Value = '  -0.17%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '42.64'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +10.37%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.328'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.05%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0696'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.11%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0989'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.17%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.123.98'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.98%  '

$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.41'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.02%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.858.05'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.02%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.679'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.56%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.69'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.41%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '35.034.22'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.67%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '70.31'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.26%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0797'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.88%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '241.27'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.73%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.17'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.15%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.75'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.56%  '

$ws.Range("E23").Value = '  -0.07%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.29'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.05%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '171.47'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.27%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.95'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +32.49%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.34%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.66'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.58%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.125'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.69%  '

$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0558'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.71%  '

$ws.Range("B31").Value = 'BinanceUSD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.01'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.20%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.00'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.43%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.99'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.66%  '

$ws.Range("E34").Value = '  +14.24%  '

$ws.Range("E35").Value = '  +22.90%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.779'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +13.30%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.25'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.47%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.08'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +14.10%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '91.59'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.17%  '

$ws.Range("E40").Value = '  +6.97%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.350.76'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.61%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '14.95'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.42%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.32'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.59%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '12.73'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +60.26%  '

$ws.Range("E45").Value = '  -1.76%  '

$ws.Range("E46").Value = '  +2.06%  '

$ws.Range("E47").Value = '  +7.61%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.33'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.06%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.044.71'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.05%  '

$ws.Range("E50").Value = '  +3.11%  '

$ws.Range("B51").Value = 'THORChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.42'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +16.89%  '
